$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "DONE ?" column header and its first value ("y")
$ws.Range("G1").Value = "DONE ?"
$ws.Range("G2").Value = "y"

# Re-apply the duplicate formulas across contiguous ranges so Excel
# collapses them into shared formulas, matching the target workbook.
$ws.Range("B5:B11").Formula = '=LEN(A5)-LEN(SUBSTITUTE(A5," ",""))+1'
$ws.Range("D5:D27").Formula = '=E5/30'

# Move the active selection to the newly added cell.
$ws.Range("G2").Select() | Out-Null
